# Insert a new date column ("11-nov") just before the "01-oct." column (column DP)
# on the "Prix Spot" sheet, shifting the existing Oct/Nov-tail columns one to the
# right (old DP:ET -> new DQ:EU). The freshly inserted column is filled with the
# header label in row 1 and with the sheet's usual "-" placeholder for the data
# rows (2-25), matching the existing convention used for other not-yet-populated
# date columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a blank column at DP (column 120); everything at/after DP shifts right.
$ws.Columns("DP:DP").Insert()

# Header row gets the new date label.
$ws.Range("DP1").Value = "11-nov"

# Data rows (2-25) get the "-" placeholder used throughout the sheet for
# dates without data yet.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 120).Value = "-"
}
